# 수정사항_애드모어 - add new "주문관리" row to the 나중에 수정할것 sheet
# (admin order-payment modal edit form) and update the current selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("나중에 수정할것")
$ws.Activate()

# Insert a new blank row above the existing row 21 (the blank separator
# row), pushing all rows from 21 downward (the old rows 22-33) down by one.
$ws.Rows.Item(21).Insert()

# Fill the new row with the new 주문관리 (order management) task entry.
$ws.Cells.Item(21, 2).Value = "관리자페이지"
$ws.Cells.Item(21, 3).Value = "주문관리"
$ws.Cells.Item(21, 4).Value = "결제내역 클릭후 모달로 나오는 폼 수정하기"

# Match the author's final cursor position/selection in the sheet.
$ws.Range("D17").Select()
